# Generate Report for Handback
# Updates the handoff/handback generated-report timestamps for the
# d6199722-e2fb-4701-ad20-c13c4e29cf31 entry across the Overview,
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the third file (row 4, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-07 07:55:23"

# zh-cn sheet: "Correspond Handoff Datetime" (col H) and
# "Correspond Handback DateTime" (col K) for row 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-07 07:55:02"
$wsZhCn.Range("K4").Value = "2016-09-07 07:56:07"

# de-de sheet: "Correspond Handback DateTime" (col K) for row 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-07 07:56:28"
